# Database.xlsx update: populate A1 on Sheet1 with the string "Cessna".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "Cessna"
